$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('E2').Value = '[Hae-Young%Lee%hylee612@snu.ac.kr%1,                   Juhee%Ahn%ahnjuhee23@gmail.com%1,                   Juhong%Park%parkjh@nhis.or.kr%1,                   Chang%Kyung Kang%zeptemiger@hanmail.net%1,                   Sung-Ho%Won%sunghow@gmail.com%1,                   Dong%Wook Kim%kimdw2269@gmail.com%1,                   Jong-Heon%Park%parkjh@nhis.or.kr%1,                   Ki-Hyun%Chung%ckhyun56@nmc.or.kr%1,                   Joon-Sung%Joh%ssabana777@gmail.com%1,                   Ji Hwan%Bang%roundbirch@gmail.com%1,                   Cheong%Hee Kang%chkangmd@nhis.or.kr%1,                   Wook%Bum Pyun%pwb423@ewha.ac.kr%1,                   Myoung-don%Oh%NULL%1]'
$ws.Range('E3').Value = '[Payam%Peymani%NULL%1,                   Tania%Dehesh%NULL%1,                   Farnaz%Aligolighasemabadi%NULL%1,                   Mohammadamin%Sadeghdoust%NULL%1,                   Katarzyna%Kotfis%NULL%1,                   Mazaher%Ahmadi%NULL%1,                   Parvaneh%Mehrbod%NULL%1,                   Pooya%Iranpour%NULL%1,                   Sanaz%Dastghaib%NULL%1,                   Ahmad%Nasimian%NULL%1,                   Amir%Ravandi%NULL%1,                   Biniam%Kidane%NULL%1,                   Naseer%Ahmed%NULL%1,                   Pawan%Sharma%NULL%1,                   Shahla%Shojaei%NULL%1,                   Kamran%Bagheri Lankarani%NULL%1,                   Andrzej%Madej%NULL%1,                   Nima%Rezaei%NULL%1,                   Tayyebeh%Madrakian%NULL%1,                   Marek J.%Los%NULL%1,                   Hagar Ibrahim%Labouta%NULL%1,                   Pooneh%Mokarram%NULL%1,                   Saeid%Ghavami%saeid.ghavami@umanitoba.ca%1]'
$ws.Range('E4').Value = '[Yongzhen%Fan%NULL%1,                   Tao%Guo%NULL%1,                   Feifei%Yan%NULL%1,                   Ming%Gong%NULL%1,                   Xin A.%Zhang%NULL%1,                   Chenze%Li%NULL%1,                   Tao%He%NULL%0,                   Huimin%Luo%NULL%1,                   Lin%Zhang%NULL%0,                   Ming%Chen%NULL%0,                   Xiaoyan%Wu%NULL%0,                   Hairong%Wang%NULL%0,                   Ke-Qiong%Deng%NULL%1,                   Jiao%Bai%NULL%1,                   Lin%Cai%NULL%1,                   Zhibing%Lu%NULL%0]'
$ws.Range('E5').Value = '[Jawad Haider%Butt%NULL%1,                   Thomas Alexander%Gerds%NULL%2,                   Thomas Alexander%Gerds%NULL%0,                   Morten%Schou%NULL%1,                   Kristian%Kragholm%NULL%1,                   Matthew%Phelps%NULL%1,                   Eva%Havers-Borgersen%NULL%1,                   Adelina%Yafasova%NULL%1,                   Gunnar Hilmar%Gislason%NULL%1,                   Christian%Torp-Pedersen%NULL%1,                   Lars%Køber%NULL%1,                   Emil Loldrup%Fosbøl%NULL%1]'
$ws.Range('E6').Value = '[Lluís%Masana%luis.masana@urv.cat%1,                   Eudald%Correig%NULL%1,                   Cèlia%Rodríguez-Borjabad%NULL%1,                   Eva%Anoro%NULL%1,                   Juan Antonio%Arroyo%NULL%1,                   Carlos%Jericó%NULL%1,                   Angels%Pedragosa%NULL%1,                   Marcel·la%Miret%NULL%1,                   Silvia%Näf%NULL%1,                   Anna%Pardo%NULL%1,                   Verónica%Perea%NULL%1,                   Rosa%Pérez-Bernalte%NULL%1,                   Núria%Plana%NULL%1,                   Rafael%Ramírez-Montesinos%NULL%1,                   Meritxell%Royuela%NULL%1,                   Cristina%Soler%NULL%1,                   Maria%Urquizu-Padilla%NULL%1,                   Alberto%Zamora%NULL%1,                   Juan%Pedro-Botet%NULL%1,                   NULL%NULL%NULL%0]'
$ws.Range('E7').Value = '[Omar%Saeed%osaeed@montefiore.org%1,                   Francesco%Castagna%NULL%2,                   Francesco%Castagna%NULL%0,                   Ilir%Agalliu%NULL%2,                   Ilir%Agalliu%NULL%0,                   Xiaonan%Xue%NULL%2,                   Xiaonan%Xue%NULL%0,                   Snehal R.%Patel%NULL%1,                   Yogita%Rochlani%NULL%1,                   Rachna%Kataria%NULL%1,                   Sasa%Vukelic%NULL%1,                   Daniel B.%Sims%NULL%1,                   Chikezie%Alvarez%NULL%1,                   Mercedes%Rivas‐Lasarte%NULL%2,                   Mercedes%Rivas‐Lasarte%NULL%0,                   Mario J.%Garcia%NULL%2,                   Mario J.%Garcia%NULL%0,                   Ulrich P.%Jorde%NULL%1]'
$ws.Range('E8').Value = '[Albano%G.%coreGivesNoEmail%0,                 Antonelli%M.%coreGivesNoEmail%1,                 Bellani%G.%coreGivesNoEmail%1,                 Bonanomi%E.%coreGivesNoEmail%1,                 Cabrini%L.%coreGivesNoEmail%1,                 Carlesso%E.%coreGivesNoEmail%1,                 Castelli%G.%coreGivesNoEmail%1,                 Cattaneo%S.%coreGivesNoEmail%1,                 Cecconi%M.%coreGivesNoEmail%1,                 Cereda%D.%coreGivesNoEmail%1,                 Colombo%S.%coreGivesNoEmail%1,                 Coluccello%A.%coreGivesNoEmail%1,                 Crescini%G.%coreGivesNoEmail%1,                 Forastieri%Molinari A.%coreGivesNoEmail%1,                 Foti%G.%coreGivesNoEmail%1,                 Fumagalli%R.%coreGivesNoEmail%1,                 Grasselli%G.%coreGivesNoEmail%1,                 Greco%M.%coreGivesNoEmail%1,                 Iotti%G. A.%coreGivesNoEmail%1,                 Langer%T.%coreGivesNoEmail%1,                 Latronico%N.%coreGivesNoEmail%1,                 Lorini%F. L.%coreGivesNoEmail%1,                 Mojoli%F.%coreGivesNoEmail%1,                 Natalini%G.%coreGivesNoEmail%1,                 Pesenti%A.%coreGivesNoEmail%1,                 Pessina%C. M.%coreGivesNoEmail%1,                 Ranieri%V. M.%coreGivesNoEmail%1,                 Rech%R.%coreGivesNoEmail%1,                 Rosano%A.%coreGivesNoEmail%1,                 Scudeller%L.%coreGivesNoEmail%1,                 Storti%E.%coreGivesNoEmail%1,                 Thompson%B. T.%coreGivesNoEmail%1,                 Tirani%M.%coreGivesNoEmail%1,                 Villani%P. G.%coreGivesNoEmail%1,                 Zanella%A.%coreGivesNoEmail%1]'
$ws.Range('E9').Value = '[Guillermo%Rodriguez-Nava%Guillermo.RodriguezNava@amitahealth.org%1,                   Daniela Patricia%Trelles-Garcia%NULL%2,                   Daniela Patricia%Trelles-Garcia%NULL%0,                   Maria Adriana%Yanez-Bello%NULL%1,                   Chul Won%Chung%NULL%1,                   Valeria Patricia%Trelles-Garcia%NULL%1,                   Harvey J.%Friedman%NULL%1]'
$ws.Range('E10').Value = '[Xiao-Jing%Zhang%NULL%1,                   Juan-Juan%Qin%NULL%1,                   Xu%Cheng%NULL%1,                   Lijun%Shen%NULL%1,                   Yan-Ci%Zhao%NULL%1,                   Yufeng%Yuan%NULL%1,                   Fang%Lei%NULL%1,                   Ming-Ming%Chen%NULL%1,                   Huilin%Yang%NULL%1,                   Liangjie%Bai%NULL%1,                   Xiaohui%Song%NULL%1,                   Lijin%Lin%NULL%1,                   Meng%Xia%NULL%1,                   Feng%Zhou%NULL%1,                   Jianghua%Zhou%NULL%1,                   Zhi-Gang%She%NULL%1,                   Lihua%Zhu%NULL%1,                   Xinliang%Ma%NULL%1,                   Qingbo%Xu%NULL%1,                   Ping%Ye%NULL%1,                   Guohua%Chen%NULL%1,                   Liming%Liu%NULL%1,                   Weiming%Mao%NULL%1,                   Youqin%Yan%NULL%1,                   Bing%Xiao%NULL%1,                   Zhigang%Lu%NULL%1,                   Gang%Peng%NULL%1,                   Mingyu%Liu%NULL%1,                   Jun%Yang%NULL%1,                   Luyu%Yang%NULL%1,                   Changjiang%Zhang%NULL%1,                   Haofeng%Lu%NULL%1,                   Xigang%Xia%NULL%1,                   Daihong%Wang%NULL%1,                   Xiaofeng%Liao%NULL%1,                   Xiang%Wei%NULL%1,                   Bing-Hong%Zhang%NULL%1,                   Xin%Zhang%NULL%2,                   Juan%Yang%NULL%1,                   Guang-Nian%Zhao%NULL%1,                   Peng%Zhang%NULL%0,                   Peter P.%Liu%NULL%1,                   Rohit%Loomba%NULL%1,                   Yan-Xiao%Ji%NULL%1,                   Jiahong%Xia%jiahong.xia@hust.edu.cn%1,                   Yibin%Wang%yibinwang@mednet.ucla.edu%1,                   Jingjing%Cai%caijingjing83@hotmail.com%1,                   Jiao%Guo%guoj@gdpu.edu.cn%1,                   Hongliang%Li%lihl@whu.edu.cn%1]'
$ws.Range('E11').Value = '[Shireen R.%Chacko%NULL%1,                   Robert%DeJoy%NULL%1,                   Kevin Bryan%Lo%NULL%0,                   Jeri%Albano%NULL%0,                   Eric%Peterson%NULL%0,                   Ruchika%Bhargav%NULL%0,                   Fahad%Gu%NULL%1,                   Grace%Salacup%NULL%1,                   Jerald%Pelayo%NULL%0,                   Zurab%Azmaiparashvili%NULL%0,                   Janani%Rangaswami%NULL%0,                   Gabriel%Patarroyo-Aponte%NULL%1,                   Sadia%Benzaquen%NULL%0,                   Ena%Gupta%NULL%1]'
$ws.Range('E12').Value = '[Christopher J.%Nicholson%NULL%1,                   Luke%Wooster%NULL%1,                   Haakon H.%Sigurslid%NULL%1,                   Rebecca H.%Li%NULL%1,                   Wanlin%Jiang%NULL%1,                   Wenjie%Tian%NULL%1,                   Christian L.%Lino Cardenas%NULL%1,                   Rajeev%Malhotra%rmalhotra@mgh.harvard.edu%1]'
$ws.Range('E13').Value = '[Aakriti%Gupta%ag3786@cumc.columbia.edu%1,                   Mahesh V.%Madhavan%NULL%2,                   Mahesh V.%Madhavan%NULL%0,                   Timothy J.%Poterucha%NULL%2,                   Timothy J.%Poterucha%NULL%0,                   Ersilia M.%DeFilippis%NULL%1,                   Jessica A.%Hennessey%NULL%1,                   Bjorn%Redfors%NULL%1,                   Christina%Eckhardt%NULL%1,                   Behnood%Bikdeli%NULL%1,                   Jonathan%Platt%NULL%1,                   Ani%Nalbandian%NULL%1,                   Pierre%Elias%NULL%2,                   Pierre%Elias%NULL%0,                   Matthew J.%Cummings%NULL%2,                   Matthew J.%Cummings%NULL%0,                   Shayan N.%Nouri%NULL%1,                   Matthew%Lawlor%NULL%1,                   Lauren S.%Ranard%NULL%1,                   Jianhua%Li%NULL%2,                   Jianhua%Li%NULL%0,                   Claudia%Boyle%NULL%1,                   Raymond%Givens%NULL%1,                   Daniel%Brodie%NULL%0,                   Harlan M.%Krumholz%NULL%2,                   Harlan M.%Krumholz%NULL%0,                   Gregg W.%Stone%NULL%2,                   Gregg W.%Stone%NULL%0,                   Sanjum S.%Sethi%NULL%1,                   Daniel%Burkhoff%NULL%1,                   Nir%Uriel%NULL%1,                   Allan%Schwartz%NULL%0,                   Martin B.%Leon%NULL%1,                   Ajay J.%Kirtane%NULL%1,                   Elaine Y.%Wan%eyw2003@cumc.columbia.edu%1,                   Sahil A.%Parikh%NULL%2,                   Sahil A.%Parikh%NULL%0]'
$ws.Range('E14').Value = '[Matthieu%Wargny%NULL%1,                   Louis%Potier%NULL%2,                   Louis%Potier%NULL%0,                   Pierre%Gourdy%NULL%2,                   Pierre%Gourdy%NULL%0,                   Matthieu%Pichelin%NULL%2,                   Matthieu%Pichelin%NULL%0,                   Coralie%Amadou%NULL%2,                   Coralie%Amadou%NULL%0,                   Pierre-Yves%Benhamou%NULL%2,                   Pierre-Yves%Benhamou%NULL%0,                   Jean-Baptiste%Bonnet%NULL%1,                   Lyse%Bordier%NULL%1,                   Olivier%Bourron%NULL%1,                   Claude%Chaumeil%NULL%2,                   Claude%Chaumeil%NULL%0,                   Nicolas%Chevalier%NULL%1,                   Patrice%Darmon%NULL%2,                   Patrice%Darmon%NULL%0,                   Blandine%Delenne%NULL%2,                   Blandine%Delenne%NULL%0,                   Delphine%Demarsy%NULL%1,                   Marie%Dumas%NULL%1,                   Olivier%Dupuy%NULL%1,                   Anna%Flaus-Furmaniuk%NULL%2,                   Anna%Flaus-Furmaniuk%NULL%0,                   Jean-François%Gautier%NULL%1,                   Anne-Marie%Guedj%NULL%2,                   Anne-Marie%Guedj%NULL%0,                   Nathalie%Jeandidier%NULL%1,                   Etienne%Larger%NULL%1,                   Jean-Philippe%Le Berre%NULL%1,                   Myriam%Lungo%NULL%1,                   Nathanaëlle%Montanier%NULL%1,                   Philippe%Moulin%NULL%1,                   Françoise%Plat%NULL%1,                   Vincent%Rigalleau%NULL%1,                   René%Robert%NULL%1,                   Dominique%Seret-Bégué%NULL%2,                   Dominique%Seret-Bégué%NULL%0,                   Pierre%Sérusclat%NULL%1,                   Sarra%Smati%NULL%1,                   Jean-François%Thébaut%NULL%1,                   Blandine%Tramunt%NULL%2,                   Blandine%Tramunt%NULL%0,                   Camille%Vatier%NULL%1,                   Fritz-Line%Velayoudom%NULL%2,                   Fritz-Line%Velayoudom%NULL%0,                   Bruno%Vergès%NULL%1,                   Patrice%Winiszewski%NULL%1,                   Audrey%Zabulon%NULL%1,                   Pierre-Antoine%Gourraud%NULL%1,                   Ronan%Roussel%NULL%2,                   Ronan%Roussel%NULL%0,                   Bertrand%Cariou%bertrand.cariou@univ-nantes.fr%2,                   Bertrand%Cariou%bertrand.cariou@univ-nantes.fr%0,                   Samy%Hadjadj%samy.hadjadj@univ-nantes.fr%1,                   NULL%NULL%NULL%0,                   NULL%NULL%NULL%0]'
$ws.Range('E15').Value = '[Tak Kyu%Oh%NULL%0,                   In-Ae%Song%NULL%1,                   Young-Tae%Jeon%NULL%1,                   Philip P.%Foster%NULL%2,                   Philip P.%Foster%NULL%0]'
$ws.Range('E16').Value = '[Gianfranco%Mitacchione%NULL%1,                   Marco%Schiavone%NULL%1,                   Antonio%Curnis%NULL%1,                   Marcello%Arca%NULL%1,                   Spinello%Antinori%NULL%1,                   Alessio%Gasperetti%NULL%1,                   Giosuè%Mascioli%NULL%1,                   Paolo%Severino%NULL%1,                   Federica%Sabato%NULL%1,                   Maria M.%Caracciolo%NULL%1,                   Gianmarco%Arabia%NULL%1,                   Laura%D''Erasmo%NULL%1,                   Maurizio%Viecca%NULL%1,                   Massimo%Mancone%NULL%1,                   Massimo%Galli%NULL%1,                   Giovanni B.%Forleo%NULL%1]'
$ws.Range('E17').Value = '[Ning%Rosenthal%NULL%1,                   Zhun%Cao%NULL%1,                   Jake%Gundrum%NULL%1,                   Jim%Sianis%NULL%1,                   Stella%Safo%NULL%1]'
$ws.Range('E18').Value = '[Maurizio%Bifulco%maubiful@unina.it%1,                   Michele%Ciccarelli%NULL%0,                   Michele%Ciccarelli%NULL%0,                   Dario%Bruzzese%NULL%1,                   Andrea%Dipasquale%NULL%1,                   Andrea G.%Lania%NULL%1,                   Gherardo%Mazziotti%NULL%1,                   Patrizia%Gazzerro%pgazzerro@unisa.it%1]'
$ws.Range('E19').Value = '[Peter J.%Mallow%NULL%1,                   Kathy W.%Belk%NULL%1,                   Michael%Topmiller%NULL%1,                   Edmond A.%Hooker%NULL%1]'
$ws.Range('E20').Value = '[Sophia L.%Song%NULL%1,                   Sarah B.%Hays%NULL%1,                   Constance E.%Panton%NULL%2,                   Constance E.%Panton%NULL%0,                   Evangelia K.%Mylona%NULL%1,                   Markos%Kalligeros%NULL%1,                   Fadi%Shehadeh%NULL%1,                   Eleftherios%Mylonakis%NULL%1]'
$ws.Range('E21').Value = '[Lori B.%Daniels%NULL%1,                   Amy M.%Sitapati%NULL%1,                   Jing%Zhang%NULL%1,                   Jingjing%Zou%NULL%1,                   Quan M.%Bui%NULL%1,                   Junting%Ren%NULL%1,                   Christopher A.%Longhurst%NULL%1,                   Michael H.%Criqui%NULL%1,                   Karen%Messer%NULL%1]'
$ws.Range('E22').Value = '[Anton%De Spiegeleer%NULL%1,                   Antoon%Bronselaer%NULL%1,                   James T.%Teo%NULL%1,                   Geert%Byttebier%NULL%1,                   Guy%De Tré%NULL%1,                   Luc%Belmans%NULL%1,                   Richard%Dobson%NULL%1,                   Evelien%Wynendaele%NULL%1,                   Christophe%Van De Wiele%NULL%1,                   Filip%Vandaele%NULL%1,                   Diemer%Van Dijck%NULL%1,                   Dan%Bean%NULL%1,                   David%Fedson%NULL%1,                   Bart%De Spiegeleer%NULL%1]'
